$d = $word.ActiveDocument

# WdColorIndex constants used below:
#   4 = wdBrightGreen  -> serialises as <w:highlight w:val="green"/>
#   7 = wdYellow       -> serialises as <w:highlight w:val="yellow"/>
$wdBrightGreen = 4
$wdYellow = 7

function Get-ParagraphContaining($doc, [string]$needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $para = $doc.Paragraphs.Item($i)
        if ($para.Range.Text.Contains($needle)) {
            return $para
        }
    }
    return $null
}

# 1) "Identificar quais regiões tem mais chamados" had no highlight at
#    all before; give the paragraph (mark + run) a green highlight.
#    Using Range.Font (not Range alone) so the highlight also lands on
#    the paragraph mark's rPr (w:pPr/w:rPr), matching the diff.
$pRegioes = Get-ParagraphContaining $d "Identificar quais regiões tem mais chamados"
$pRegioes.Range.Font.HighlightColorIndex = $wdBrightGreen

# 2) In the "ambulâncias necessário" paragraph, insert a new clause
#    " (Manhã/Tarde/Noite)" right after "...horário do dia" and before
#    " por cidade/bairro...".
$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute(
    "necessário para um certo horário do dia por cidade/bairro",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "necessário para um certo horário do dia (Manhã/Tarde/Noite) por cidade/bairro",
    2
) | Out-Null

#    The whole paragraph (incl. the "rua ," run and trailing run, and
#    the paragraph mark) changes highlight from green to yellow.
$pAmbulancias = Get-ParagraphContaining $d "Identificar o número médio de ambulâncias necessário"
$pAmbulancias.Range.Font.HighlightColorIndex = $wdYellow

# 3) "Plotar um mapa indicando problemas de saúdes por regiões." goes
#    from green to yellow as well.
$pMapa = Get-ParagraphContaining $d "Plotar um mapa indicando problemas de saúdes por regiões."
$pMapa.Range.Font.HighlightColorIndex = $wdYellow
